# Insert a new data row at row 19 (pushing the existing rows 19-59 down to
# 20-60, same as Excel's "Insert Sheet Rows" above the selection), then fill
# the newly-inserted row with the Vega Modelo de Temuco / Bruselas (repollito)
# record for 2021-10-08 (serial 44477) with a volume of 40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(19).Insert()

$ws.Range("A19").Value = 10
$ws.Range("B19").Value = "Vega Modelo de Temuco"
$ws.Range("C19").Value = "La Araucanía"
$ws.Range("D19").Value = 44477
$ws.Range("E19").Value = 9
$ws.Range("F19").Value = 100112035
$ws.Range("G19").Value = "Bruselas (repollito)"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 40
$ws.Range("K19").Value = 25000
$ws.Range("L19").Value = 25000
$ws.Range("M19").Value = 25000
$ws.Range("N19").Value = "$/malla 10 kilos"
$ws.Range("O19").Value = "Provincia de Quillota"
$ws.Range("P19").Value = 2500
$ws.Range("Q19").Value = 10
$ws.Range("R19").Value = "Hortaliza"
